$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 2; B = 1.02; C = 1.031589319445112; D = 1.035522584672375; E = 1.041374200652114; F = 1.052929954337091; I = 1.035188598823858; J = 1.036724212586162; K = 1.038318947439617; L = 1.044153896628837; M = 1.055677326738416; N = 1.038196479881963 },
    @{ Row = 3; B = 1.02; C = 1.032594152160049; D = 1.036278648432915; E = 1.042293864788724; F = 1.053999318308998; I = 1.035403917833929; J = 1.037370702032726; K = 1.038884565365507; L = 1.044883889048167; M = 1.056558976903543; N = 1.038843887417695 },
    @{ Row = 4; B = 1.02; C = 1.033244448928516; D = 1.036767540901947; E = 1.042889414533445; F = 1.054691808538576; I = 1.035541451957519; J = 1.037788545945816; K = 1.039249566115976; L = 1.045356080556332; M = 1.057129404789825; N = 1.039262324717069 },
    @{ Row = 5; B = 1.02; C = 1.033517857992554; D = 1.036972990989057; E = 1.043139894344749; F = 1.05498306013602; I = 1.03559884187697; J = 1.037964092533779; K = 1.039402774260624; L = 1.045554550389553; M = 1.05736919799224; N = 1.039438120601322 },
    @{ Row = 6; B = 1.02; C = 1.033563765993563; D = 1.037007482224771; E = 1.043181957476648; F = 1.055031970089842; I = 1.035608452700735; J = 1.037993560818708; K = 1.039428484594915; L = 1.0455878720155; M = 1.057409459447044; N = 1.039467630734596 },
    @{ Row = 7; B = 1.02; C = 1.033248102138708; D = 1.036770286452867; E = 1.042892761022928; F = 1.054695699751748; I = 1.035542220492414; J = 1.037790892060473; K = 1.039251614227463; L = 1.045358732675832; M = 1.057132608976808; N = 1.039264674163478 },
    @{ Row = 8; B = 1.02; C = 1.031928885891911; D = 1.035778167977325; E = 1.041684908604092; F = 1.053291238662475; I = 1.035261737572611; J = 1.036942795355549; K = 1.038510305610685; L = 1.044400634252387; M = 1.055975296325919; N = 1.038415373063955 },
    @{ Row = 9; B = 1.02; C = 1.029605058909363; D = 1.034027433674783; E = 1.039560116821161; F = 1.050820563553192; I = 1.034753790292194; J = 1.035444704159078; K = 1.037196459673028; L = 1.042711124609655; M = 1.053935540797723; N = 1.036915154406151 },
    @{ Row = 10; B = 1.02; C = 1.028056394136447; D = 1.032858656056731; E = 1.038146050733157; F = 1.049176281095294; I = 1.03440597520213; J = 1.034443562781163; K = 1.036315515956699; L = 1.041583999203535; M = 1.052575453623614; N = 1.035912591292647 },
    @{ Row = 11; B = 1.02; C = 1.027385938637063; D = 1.032352189066475; E = 1.037534336675717; F = 1.048464967256927; I = 1.034253194710723; J = 1.034009490912659; K = 1.035932868552277; L = 1.041095761493912; M = 1.051986467977575; N = 1.035477902992301 },
    @{ Row = 12; B = 1.02; C = 1.027136920697609; D = 1.032164008632502; E = 1.037307207455477; F = 1.048200854903179; I = 1.034196118957356; J = 1.033848171905701; K = 1.035790557601079; L = 1.040914380920621; M = 1.051767684062822; N = 1.03531635489385 },
    @{ Row = 13; B = 1.02; C = 1.027190335026178; D = 1.032204376497165; E = 1.037355923424603; F = 1.048257503309895; I = 1.034208376658793; J = 1.033882779237813; K = 1.035821091859332; L = 1.040953288899861; M = 1.051814614329391; N = 1.035351011372343 },
    @{ Row = 14; B = 1.02; C = 1.027365354346773; D = 1.032336635151575; E = 1.037515560292328; F = 1.048443133578934; I = 1.034248483464884; J = 1.033996157972775; K = 1.035921108726197; L = 1.041080769074516; M = 1.051968383400234; N = 1.035464551118113 },
    @{ Row = 15; B = 1.02; C = 1.027473192057601; D = 1.032418116672147; E = 1.037613929595581; F = 1.048557519946973; I = 1.034273151368402; J = 1.034066003050696; K = 1.035982708760332; L = 1.041159310206609; M = 1.052063124496776; N = 1.035534495384056 },
    @{ Row = 16; B = 1.02; C = 1.028100892732432; D = 1.032892260698047; E = 1.038186660575308; F = 1.049223502831795; I = 1.034416068963653; J = 1.034472358704875; K = 1.036340885918392; L = 1.041616398072894; M = 1.052614541456579; N = 1.035941428109874 },
    @{ Row = 17; B = 1.02; C = 1.028494666717268; D = 1.033189577954865; E = 1.038546076867259; F = 1.049641436327734; I = 1.034505135554425; J = 1.034727102272992; K = 1.036565241994815; L = 1.041903068005353; M = 1.052960415277017; N = 1.036196533443077 },
    @{ Row = 18; B = 1.02; C = 1.028724360627731; D = 1.033362961498697; E = 1.038755774799756; F = 1.04988527460935; I = 1.034556876779738; J = 1.034875634761761; K = 1.036695989957506; L = 1.042070259850478; M = 1.053162151910607; N = 1.036345276865017 },
    @{ Row = 19; B = 1.02; C = 1.0288026823838; D = 1.033422074576174; E = 1.038827285939268; F = 1.049968428131002; I = 1.034574483601335; J = 1.034926271139549; K = 1.036740552086915; L = 1.042127264903448; M = 1.053230937906936; N = 1.036395985152269 },
    @{ Row = 20; B = 1.02; C = 1.028452417232058; D = 1.033157682404727; E = 1.038507509053184; F = 1.049596589318953; I = 1.034495601250994; J = 1.034699776397789; K = 1.036541182616526; L = 1.041872312880392; M = 1.052923306877223; N = 1.036169168761997 },
    @{ Row = 21; B = 1.02; C = 1.027313814991698; D = 1.032297689814728; E = 1.037468548762601; F = 1.048388467292312; I = 1.034236682013898; J = 1.033962773115104; K = 1.035891661179382; L = 1.041043230091883; M = 1.051923102470476; N = 1.035431118850115 },
    @{ Row = 22; B = 1.02; C = 1.026598041026923; D = 1.031756654498414; E = 1.036815826195907; F = 1.04762945918823; I = 1.03407200166572; J = 1.033498896069976; K = 1.035482247829924; L = 1.040521794819213; M = 1.05129418575137; N = 1.034966583046375 },
    @{ Row = 23; B = 1.02; C = 1.026977475855501; D = 1.032043497976542; E = 1.037161797954008; F = 1.048031768043997; I = 1.034159480625158; J = 1.033744852743281; K = 1.035699383365131; L = 1.040798232354895; M = 1.051627590858465; N = 1.035212889006368 },
    @{ Row = 24; B = 1.02; C = 1.028471507917967; D = 1.033172094743953; E = 1.038524936013662; F = 1.04961685355168; I = 1.034499910040581; J = 1.034712123954848; K = 1.036552054369698; L = 1.041886209852717; M = 1.05294007458337; N = 1.036181533854003 },
    @{ Row = 25; B = 1.02; C = 1.030205727466834; D = 1.034480330311816; E = 1.040108995611398; F = 1.051458795217254; I = 1.034886727534127; J = 1.035832424369114; K = 1.037537012725857; L = 1.054462913069983; M = 1.037303425223358 }
)

$cols = @("B","C","D","E","F","I","J","K","L","M","N")

foreach ($r in $rows) {
    $rowNum = $r.Row
    foreach ($col in $cols) {
        if ($r.ContainsKey($col)) {
            $addr = "$col$rowNum"
            $ws.Range($addr).Value = $r[$col]
        }
    }
}
